# Update the as-of date in the confidentiality footer (row 59 of the model holdings sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; unprotect so the Weight/Percent Change cells can be refreshed,
# then restore protection once the data has been rewritten.
$ws.Unprotect()

# Each entry is Row, Weight (col D), Percent Change (col E)
$updates = @(
    @(2, 0.01104733671371894, 0.005994358251057763),
    @(3, 0.01031722862544576, -0.0009708737864076999),
    @(4, 0.0106471128287448, 0.01404917210235834),
    @(5, 0.01131667536823436, 0.008536585365853666),
    @(6, 0.01065779733735368, 0.01528822055137846),
    @(7, 0.01119124368904474, 0.08253358925143961),
    @(8, 0.01126191726161387, -0.0008103727714747766),
    @(9, 0.01131311386536473, 0.0008657327246970326),
    @(10, 0.01042407371153452, 0.004356181934657366),
    @(11, 0.01107015259147748, 0.0003016136329361263),
    @(12, 0.4496709004403019, 0.001915708812260775),
    @(13, 0.01157555210807024, 0.03072899640405358),
    @(14, 0.0109772196259732, -0.00823278921220727),
    @(15, 0.01048461926031815, -0.005137786081270357),
    @(16, 0.01003831843196823, 0.02219659844335564),
    @(17, 0.01049196485998675, -0.01983663943990666),
    @(18, 0.009189790373280002, -0.06518105849582156),
    @(19, 0.008269920960234597, 0.07076239822353814),
    @(20, 0.01222708453928232, 0.01929728745676318),
    @(21, 0.01226448031941339, -0.02107154524665145),
    @(22, 0.01186670496766211, 0.005008347245408995),
    @(23, 0.01186147401032235, 0.004926108374384341),
    @(24, 0.01066336218558747, 0.1167936541070869),
    @(25, 0.01187761207020034, 0.045127436281859),
    @(26, 0.01103576182939266, 0.0308604623018276),
    @(27, 0.01028539769354848, 0.0119462419113987),
    @(28, 0.01239569944076614, 0.01125925925925952),
    @(29, 0.01003342136552249, 0.02628951747088193),
    @(30, 0.006999911296319154, 0.004754030591153313),
    @(31, 0.00537564339384067, -0.006521739130434745),
    @(32, 0.009400252933481923, -0.03793466807165435),
    @(33, 0.01067371280330232, 0.008466888418506402),
    @(34, 0.01060281663680384, -0.01183003380009662),
    @(35, 0.009545495472383829, 0.00690250215703192),
    @(36, 0.01122006960289577, -0.01414514145141454),
    @(37, 0.009870816500131165, 0.01035077630822312),
    @(38, 0.01124488882601847, 0.02558519324986386),
    @(39, 0.01352124564453241, -0.008823916765441919),
    @(40, 0.01124143862011352, 0.004514672686230403),
    @(41, 0.01158200733202144, 0.01937269372693717),
    @(42, 0.01111456008038312, 0.01598173515981727),
    @(43, 0.01103130995080563, 0.01074498567335236),
    @(44, 0.0110098296366232, 0.01069518716577544),
    @(45, 0.01157488432628219, 0.02307692307692299),
    @(46, 0.01075885498127929, 0.002048247610378029),
    @(47, 0.01055028446947687, -0.006202924235711227),
    @(48, 0.01081806496648681, 0.01234567901234551),
    @(49, 0.009975658240855759, 0.0007698229407235857),
    @(50, 0.009604037675803295, 0),
    @(51, 0.01001617033599775, -0.01594533029612744),
    @(52, 0.01013915348196449, -0.001097694840834218),
    @(53, 0.009441321513447289, -0.01002004008016033),
    @(54, 0.004151933267230357, 0.009248090068355364),
    @(55, 0.004079701537155768, 0.003928415539065977),
    @(56, 1, 0.005883046255797941)
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 4).Value = $u[1]
    $ws.Cells.Item($r, 5).Value = $u[2]
}

# Refresh the "as of" date in the confidentiality disclosure (2021-04-28 -> 2021-04-29)
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-29 for illustrative purposes only and are subject to change."

$ws.Protect()
